# Edit script: applies the diff changes to the betexplorer Colombia Primera-A 2023 workbook.
# 1) Rows 192-196 and rows 197-201 get their match-data columns (F:V) permuted
#    (the match rows were re-ordered while keeping Indice/pais/torneio/temporada/data_partida
#    in columns A:E as-is).
# 2) Ten new match rows (218-227) are appended at the end of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Step 1: permute F:V across rows 192-201 (two blocks of 5 rows each)
# ---------------------------------------------------------------------------

# capture original F:V values for rows 192-201 before overwriting anything
$orig = @{}
for ($r = 192; $r -le 201; $r++) {
    $orig[$r] = $ws.Range("F$r`:V$r").Value()
}

# mapping: new row number -> row number that currently holds the data that should end up there
$mapping = @{
    192 = 195
    193 = 194
    194 = 196
    195 = 192
    196 = 193
    197 = 200
    198 = 201
    199 = 198
    200 = 199
    201 = 197
}

foreach ($newR in $mapping.Keys) {
    $oldR = $mapping[$newR]
    $ws.Range("F$newR`:V$newR").Value = $orig[$oldR]
}

# ---------------------------------------------------------------------------
# Step 2: append new rows 218-227
# ---------------------------------------------------------------------------

# Copy formatting (styles) from the last existing data row (217) into the new rows
$srcRange = $ws.Range("A217:V217")
for ($r = 218; $r -le 227; $r++) {
    $destRange = $ws.Range("A$r`:V$r")
    $srcRange.Copy($destRange)
}

# Column D (temporada) holds text that looks numeric ("2023"); force text format
# for the new rows so Excel does not silently convert it to a number.
$ws.Range("D218:D227").NumberFormat = "@"

$newRowsData = @(
    @(217,"colombia","primera-a","2023",45263,"Deportes Tolima",0,"Aguilas",1,2.06,"30/11/2023 02:42",1.85,"02/12/2023 23:59",3.18,"30/11/2023 02:42",3.41,"02/12/2023 23:59",3.78,"30/11/2023 02:42",4.93,"02/12/2023 23:59","https://www.betexplorer.com/football/colombia/primera-a/deportes-tolima-aguilas-doradas/8Q5unWbN/"),
    @(218,"colombia","primera-a","2023",45263.09375,"Dep. Cali",0,"Junior",2,2.57,"30/11/2023 02:42",3.91,"03/12/2023 02:14",2.99,"30/11/2023 02:42",3.29,"03/12/2023 02:14",2.95,"30/11/2023 02:42",2.11,"03/12/2023 02:14","https://www.betexplorer.com/football/colombia/primera-a/dep-cali-junior/nH4qojET/"),
    @(219,"colombia","primera-a","2023",45263.89583333334,"Atl. Nacional",0,"Ind. Medellin",5,2.26,"01/12/2023 02:42",3.08,"03/12/2023 21:26",3.16,"01/12/2023 02:42",3.24,"03/12/2023 21:26",3.28,"01/12/2023 02:42",2.5,"03/12/2023 21:26","https://www.betexplorer.com/football/colombia/primera-a/atl-nacional-ind-medellin/v9dBtAMp/"),
    @(220,"colombia","primera-a","2023",45264.04166666666,"America De Cali",1,"Millonarios",0,2.06,"01/12/2023 02:42",2.39,"04/12/2023 00:55",3.2,"01/12/2023 02:42",3.21,"04/12/2023 00:50",3.76,"01/12/2023 02:42",3.3,"04/12/2023 00:55","https://www.betexplorer.com/football/colombia/primera-a/america-de-cali-millonarios/IaeFuUyi/"),
    @(221,"colombia","primera-a","2023",45267,"Millonarios",0,"Atl. Nacional",1,1.87,"04/12/2023 01:12",1.87,"06/12/2023 23:52",3.31,"04/12/2023 01:12",3.43,"06/12/2023 23:52",4.36,"04/12/2023 01:12",4.78,"06/12/2023 23:43","https://www.betexplorer.com/football/colombia/primera-a/millonarios-atl-nacional/MLZOKVR9/"),
    @(222,"colombia","primera-a","2023",45267,"Aguilas",3,"Dep. Cali",1,1.79,"03/12/2023 08:12",1.36,"06/12/2023 23:15",3.5,"03/12/2023 08:12",5.02,"06/12/2023 23:15",4.45,"03/12/2023 08:12",9.46,"06/12/2023 23:15","https://www.betexplorer.com/football/colombia/primera-a/aguilas-doradas-dep-cali/UoZSJksG/"),
    @(223,"colombia","primera-a","2023",45267.09375,"Ind. Medellin",2,"America De Cali",1,2.05,"04/12/2023 01:12",2.29,"07/12/2023 02:12",3.31,"04/12/2023 01:12",3.38,"07/12/2023 02:12",3.64,"04/12/2023 01:12",3.31,"07/12/2023 02:12","https://www.betexplorer.com/football/colombia/primera-a/ind-medellin-america-de-cali/rVzLLBC3/"),
    @(224,"colombia","primera-a","2023",45267.09375,"Junior",4,"Deportes Tolima",2,2.06,"03/12/2023 08:12",1.93,"07/12/2023 02:05",3.2,"03/12/2023 08:12",3.63,"07/12/2023 02:13",3.75,"03/12/2023 08:12",4.12,"07/12/2023 02:13","https://www.betexplorer.com/football/colombia/primera-a/junior-deportes-tolima/8xYWI9dM/"),
    @(225,"colombia","primera-a","2023",45270.91666666666,"Junior",3,"Ind. Medellin",2,1.87,"07/12/2023 17:43",1.99,"10/12/2023 21:56",3.35,"07/12/2023 17:43",3.24,"10/12/2023 21:56",4.27,"07/12/2023 17:43",4.45,"10/12/2023 21:56","https://www.betexplorer.com/football/colombia/primera-a/junior-ind-medellin/jFjvuB1g/"),
    @(226,"colombia","primera-a","2023",45274.08333333334,"Ind. Medellin",2,"Junior",1,2.16,"11/12/2023 07:43",1.74,"14/12/2023 01:57",3.12,"11/12/2023 07:43",3.65,"14/12/2023 01:57",3.57,"11/12/2023 07:43",5.34,"14/12/2023 01:57","https://www.betexplorer.com/football/colombia/primera-a/ind-medellin-junior/U5krvVGa/")
)

$nrows = $newRowsData.Count
$ncols = $newRowsData[0].Count
$arr = New-Object 'object[,]' $nrows,$ncols
for ($i = 0; $i -lt $nrows; $i++) {
    for ($j = 0; $j -lt $ncols; $j++) {
        $arr[$i,$j] = $newRowsData[$i][$j]
    }
}

$ws.Range("A218:V227").Value = $arr

# ---------------------------------------------------------------------------
# Step 3: make sure the sheet dimension reflects the new extent
# ---------------------------------------------------------------------------
Write-Output "edit complete"
